# Penalty method -> steepest descent method: refreshed results.
# Updates the per-epsilon convergence values (columns B-D), the
# iteration-count columns (E-F), and the epsilon values in column A
# that changed because the new solver took a different path / used a
# different epsilon schedule. Also narrows columns E and F since the
# iteration counts are now much smaller (single/double digits instead
# of triple/quadruple digits).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (epsilon = 0.5, unchanged)
$ws.Range("B1").Value = 0.41036011658012789
$ws.Range("C1").Value = 0.41036011658012789
$ws.Range("D1").Value = 0.41036011658012789
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 19

# Row 2
$ws.Range("A2").Value = 0.14999999999999999
$ws.Range("B2").Value = 0.40745893220454099
$ws.Range("C2").Value = 0.40745893220454099
$ws.Range("D2").Value = 0.40745893220454099
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 15

# Row 3
$ws.Range("A3").Value = 0.044999999999999998
$ws.Range("B3").Value = 0.40856708602838504
$ws.Range("C3").Value = 0.40856708602838504
$ws.Range("D3").Value = 0.40856708602838504
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 39

# Row 4
$ws.Range("A4").Value = 0.0135
$ws.Range("B4").Value = 0.40840540856439328
$ws.Range("C4").Value = 0.40840540856439328
$ws.Range("D4").Value = 0.40840540856439328
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 49

# Row 5
$ws.Range("A5").Value = 0.0040499999999999998
$ws.Range("B5").Value = 0.40830548639643149
$ws.Range("C5").Value = 0.40830548639643149
$ws.Range("D5").Value = 0.40830548639643149
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 49

# The iteration-count columns got a lot narrower now that the values
# are 1-2 digits instead of 3-4 digits.
$ws.Columns.Item(5).ColumnWidth = 1.3
$ws.Columns.Item(6).ColumnWidth = 2.25
